$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.027239793340086
$ws.Range("D2").Value = 1.032067801183329
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.025769974936533
$ws.Range("I2").Value = 1.033403940087093
$ws.Range("J2").Value = 1.032398770331024
$ws.Range("K2").Value = 1.034874100395477
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.028594574845461
$ws.Range("N2").Value = 1.033864895002695
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.028112611192807
$ws.Range("D3").Value = 1.032720676740583
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.027275593639262
$ws.Range("I3").Value = 1.033592032695794
$ws.Range("J3").Value = 1.032912082982618
$ws.Range("K3").Value = 1.035336088330424
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.029905673537739
$ws.Range("N3").Value = 1.034378936617125
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.028677286501425
$ws.Range("D4").Value = 1.033143012491502
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.028249758446725
$ws.Range("I4").Value = 1.033712509755861
$ws.Range("J4").Value = 1.033243504317966
$ws.Range("K4").Value = 1.035634223753689
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.030753482017496
$ws.Range("N4").Value = 1.034710828608784
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028914652776369
$ws.Range("D5").Value = 1.03332053341386
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.028659284792242
$ws.Range("I5").Value = 1.033762863349811
$ws.Range("J5").Value = 1.033382659882959
$ws.Range("K5").Value = 1.035759367684697
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.031109770439057
$ws.Range("N5").Value = 1.034850181790641
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028954506250815
$ws.Range("D6").Value = 1.033350338235661
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.028728045462978
$ws.Range("I6").Value = 1.033771300638585
$ws.Range("J6").Value = 1.03340601451051
$ws.Range("K6").Value = 1.035780368622366
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.031169585304779
$ws.Range("N6").Value = 1.034873569584442
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.028680458292711
$ws.Range("D7").Value = 1.033145384649685
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.028255230601465
$ws.Range("I7").Value = 1.033713183742241
$ws.Range("J7").Value = 1.033245364405122
$ws.Range("K7").Value = 1.035635896688599
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.030758243267925
$ws.Range("N7").Value = 1.034712691337477
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.027534786025475
$ws.Range("D8").Value = 1.032288467639116
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.026278823272501
$ws.Range("I8").Value = 1.033467761749512
$ws.Range("J8").Value = 1.03257239707061
$ws.Range("K8").Value = 1.035030397203136
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.029037785044448
$ws.Range("N8").Value = 1.03403876831216
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.025515236114998
$ws.Range("D9").Value = 1.030777593288499
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.022795394449915
$ws.Range("I9").Value = 1.033025872071105
$ws.Range("J9").Value = 1.031380983763311
$ws.Range("K9").Value = 1.033957301875575
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.026001657715081
$ws.Range("N9").Value = 1.032845663061307
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024168386948189
$ws.Range("D10").Value = 1.029769792084236
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.020472323463416
$ws.Range("I10").Value = 1.032724950046439
$ws.Range("J10").Value = 1.030582973258104
$ws.Range("K10").Value = 1.03323779914604
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.023974333680719
$ws.Range("N10").Value = 1.03204651928965
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023585072567535
$ws.Range("D11").Value = 1.029333279619038
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.01946615841941
$ws.Range("I11").Value = 1.032593147639104
$ws.Range("J11").Value = 1.030236539818349
$ws.Range("K11").Value = 1.032925275053062
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.023095657726011
$ws.Range("N11").Value = 1.031699593874674
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023368385574309
$ws.Range("D12").Value = 1.029171120770024
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.019092379327999
$ws.Range("I12").Value = 1.032543964869332
$ws.Range("J12").Value = 1.030107725252209
$ws.Range("K12").Value = 1.032809043372926
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.022769148736345
$ws.Range("N12").Value = 1.031570596377074
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023414866488223
$ws.Range("D13").Value = 1.029205905230645
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.019172558319037
$ws.Range("I13").Value = 1.032554524943409
$ws.Range("J13").Value = 1.030135362468672
$ws.Range("K13").Value = 1.032833982082231
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.02283919196454
$ws.Range("N13").Value = 1.031598272841556
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023567161521874
$ws.Range("D14").Value = 1.029319875896096
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.019435262664923
$ws.Range("I14").Value = 1.032589086775606
$ws.Range("J14").Value = 1.030225894695457
$ws.Range("K14").Value = 1.032915670290297
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.023068671062631
$ws.Range("N14").Value = 1.031688933634486
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023660993123407
$ws.Range("D15").Value = 1.029390094524544
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.019597117407594
$ws.Range("I15").Value = 1.032610351583299
$ws.Range("J15").Value = 1.030281656860344
$ws.Range("K15").Value = 1.032965981706283
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.023210043391474
$ws.Range("N15").Value = 1.031744774988043
# Row 16
$ws.Range("B16").Value = 1.019999999999999
$ws.Range("C16").Value = 1.024207097073332
$ws.Range("D16").Value = 1.029798759329872
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.02053909337419
$ws.Range("I16").Value = 1.032733665714106
$ws.Range("J16").Value = 1.030605946152296
$ws.Range("K16").Value = 1.033258519815755
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.024032630592916
$ws.Range("N16").Value = 1.032069524807987
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024549621537389
$ws.Range("D17").Value = 1.030055070102363
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.02112989650787
$ws.Range("I17").Value = 1.032810615535791
$ws.Range("J17").Value = 1.030809126007173
$ws.Range("K17").Value = 1.033441760381405
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.024548391926353
$ws.Range("N17").Value = 1.032272993201563
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024749398707359
$ws.Range("D18").Value = 1.030204559464565
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.021474477268564
$ws.Range("I18").Value = 1.032855354219148
$ws.Range("J18").Value = 1.030927551550502
$ws.Range("K18").Value = 1.033548547454906
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.02484914689877
$ws.Range("N18").Value = 1.032391586922747
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024817515601469
$ws.Range("D19").Value = 1.030255529358635
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.02159196637978
$ws.Range("I19").Value = 1.032870584384154
$ws.Range("J19").Value = 1.030967917045506
$ws.Range("K19").Value = 1.033584943131533
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.024951683165445
$ws.Range("N19").Value = 1.032432009741383
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024512873099495
$ws.Range("D20").Value = 1.030027571658249
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.021066511518348
$ws.Range("I20").Value = 1.032802374529843
$ws.Range("J20").Value = 1.030787335609918
$ws.Range("K20").Value = 1.033422110129811
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.024493063913137
$ws.Range("N20").Value = 1.032251171859444
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023522314929384
$ws.Range("D21").Value = 1.029286314887306
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.019357904058383
$ws.Range("I21").Value = 1.032578915401621
$ws.Range("J21").Value = 1.030199238904478
$ws.Range("K21").Value = 1.03289161918721
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.023001098803983
$ws.Range("N21").Value = 1.031662239989226
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.0228994071677
$ws.Range("D22").Value = 1.028820149493308
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.018283373967323
$ws.Range("I22").Value = 1.032437113202197
$ws.Range("J22").Value = 1.029828705204019
$ws.Range("K22").Value = 1.032557231985072
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.022062288845682
$ws.Range("N22").Value = 1.031291180088412
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02322963228771
$ws.Range("D23").Value = 1.029067282636017
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.018853029271652
$ws.Range("I23").Value = 1.032512408874064
$ws.Range("J23").Value = 1.030025205528712
$ws.Range("K23").Value = 1.032734577168431
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.022560042658781
$ws.Range("N23").Value = 1.031487959466105
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024529478172318
$ws.Range("D24").Value = 1.030039997059702
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.02109515254179
$ws.Range("I24").Value = 1.0328060987332
$ws.Range("J24").Value = 1.030797182017552
$ws.Range("K24").Value = 1.033430989523186
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.024518064506444
$ws.Range("N24").Value = 1.032261032250106
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.026037423762437
$ws.Range("D25").Value = 1.031168290509904
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.023696061872705
$ws.Range("I25").Value = 1.033141227220406
$ws.Range("J25").Value = 1.031689651136755
$ws.Range("K25").Value = 1.034235447440588
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.026787122659921
$ws.Range("N25").Value = 1.033154768777826

Write-Host "Applied 380 kV case updates"